$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.395.93"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "2.607.19"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E11").Value = "  -2.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.372"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "3.064.72"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.91%  "
$ws.Range("D15").Value = "60.382.61"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").Value = "2.612.86"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.18%  "
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  +3.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("E27").Value = "  +4.07%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.56%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0801"
$ws.Range("E29").Value = "  +1.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +11.14%  "
$ws.Range("E35").Value = "  +1.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.988"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.16%  "
$ws.Range("E37").Value = "  +5.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.43%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "311.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.842"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "135.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.88%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0553"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0243"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("E50").Value = "  +2.97%  "
$ws.Range("E51").Value = "  +0.55%  "
